# S0155_DeleteCase.xlsx - replace the four caseID values in column A
# with a new batch of case IDs, then select A2:A4 (active cell A2) to
# match the saved UI selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The caseID values are zero-padded numeric-looking strings that must be
# stored as text (shared strings), not numbers. Temporarily force a text
# number format so Excel doesn't strip the leading zeros, then clear the
# format back off again so the cells keep their original (default) style.
$ws.Range("A2:A4").NumberFormat = "@"
$ws.Range("A2").Value = "00001324"
$ws.Range("A3").Value = "00001325"
$ws.Range("A4").Value = "00001327"
$ws.Range("A2:A4").ClearFormats()

# Match the workbook's saved selection/active cell (A2 active, A2:A4 selected).
$ws.Range("A2:A4").Select()
